$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2-33, columns A-D
$data = New-Object 'object[,]' 32,4
$data[0,0] = 16
$data[0,1] = '{ankle, hip} (distance) - |temporal| std_min'
$data[0,2] = 0.957
$data[0,3] = 0.001
$data[1,0] = 18
$data[1,1] = '{ankle, hip} (distance) {diff} - |spectral| rel_pwr_4_to_6_min'
$data[1,2] = 1
$data[1,3] = 0.244
$data[2,0] = 27
$data[2,1] = '{ankle, hip} (distance) {diff} - |temporal| std_min'
$data[2,2] = 1
$data[2,3] = 0.001
$data[3,0] = 32
$data[3,1] = '{elbow, shoulder, hip} (angle) - |spectral| rel_pwr_2_to_4_min'
$data[3,2] = 0.826
$data[3,3] = 0
$data[4,0] = 34
$data[4,1] = '{elbow, shoulder, hip} (angle) - |spectral| rel_pwr_6_to_128_min'
$data[4,2] = 0.913
$data[4,3] = 0.003
$data[5,0] = 44
$data[5,1] = '{elbow, shoulder, hip} (angle) {diff} - |spectral| rel_pwr_0.5_to_1_max'
$data[5,2] = 0.87
$data[5,3] = 0.002
$data[6,0] = 49
$data[6,1] = '{elbow, shoulder, hip} (angle) {diff} - |spectral| rel_pwr_4_to_6_max'
$data[6,2] = 0.957
$data[6,3] = 0.003
$data[7,0] = 50
$data[7,1] = '{elbow, shoulder, hip} (angle) {diff} - |spectral| rel_pwr_4_to_6_min'
$data[7,2] = 1
$data[7,3] = 0.006
$data[8,0] = 54
$data[8,1] = '{elbow, shoulder, hip} (angle) {diff} - |temporal| median_max'
$data[8,2] = 0.913
$data[8,3] = 0.466
$data[9,0] = 59
$data[9,1] = '{left_ankle, right_ankle} (x_displacement) - |spectral| rel_pwr_0.5_to_1'
$data[9,2] = 0.957
$data[9,3] = 0.001
$data[10,0] = 60
$data[10,1] = '{left_ankle, right_ankle} (x_displacement) - |spectral| rel_pwr_1_to_2'
$data[10,2] = 0.957
$data[10,3] = 0.019
$data[11,0] = 61
$data[11,1] = '{left_ankle, right_ankle} (x_displacement) - |spectral| rel_pwr_2_to_4'
$data[11,2] = 0.87
$data[11,3] = 0.059
$data[12,0] = 64
$data[12,1] = '{left_ankle, right_ankle} (x_displacement) {diff} - |spectral| rel_pwr_1_to_2'
$data[12,2] = 1
$data[12,3] = 0.052
$data[13,0] = 69
$data[13,1] = '{left_ankle, right_ankle} (y_displacement) - |spectral| rel_pwr_2_to_4'
$data[13,2] = 1
$data[13,3] = 0.053
$data[14,0] = 71
$data[14,1] = '{left_ankle, right_ankle} (y_displacement) - |spectral| rel_pwr_6_to_128'
$data[14,2] = 1
$data[14,3] = 0.006
$data[15,0] = 81
$data[15,1] = '{left_hip, right_hip} (y_displacement) - |spectral| half_pwr_freq'
$data[15,2] = 1
$data[15,3] = 0.012
$data[16,0] = 87
$data[16,1] = '{left_hip, right_hip} (y_displacement) - |temporal| median'
$data[16,2] = 0.87
$data[16,3] = 0.124
$data[17,0] = 90
$data[17,1] = '{left_hip, right_hip} (y_displacement) {diff} - |spectral| rel_pwr_0.5_to_1'
$data[17,2] = 0.957
$data[17,3] = 0.075
$data[18,0] = 91
$data[18,1] = '{left_hip, right_hip} (y_displacement) {diff} - |spectral| rel_pwr_2_to_4'
$data[18,2] = 0.783
$data[18,3] = 0.026
$data[19,0] = 94
$data[19,1] = '{left_knee, right_knee} (x_displacement) - |spectral| entropy'
$data[19,2] = 0.826
$data[19,3] = 0.078
$data[20,0] = 96
$data[20,1] = '{left_knee, right_knee} (x_displacement) - |spectral| rel_pwr_0.5_to_1'
$data[20,2] = 0.957
$data[20,3] = 0.004
$data[21,0] = 104
$data[21,1] = '{left_knee, right_knee} (x_displacement) {diff} - |temporal| kurtosis'
$data[21,2] = 1
$data[21,3] = 0.058
$data[22,0] = 105
$data[22,1] = '{left_knee, right_knee} (x_displacement) {diff} - |temporal| median'
$data[22,2] = 0.913
$data[22,3] = 0.457
$data[23,0] = 114
$data[23,1] = '{left_shoulder, right_shoulder} (y_displacement) {diff} - |spectral| rel_pwr_0.5_to_1'
$data[23,2] = 0.826
$data[23,3] = 0.015
$data[24,0] = 116
$data[24,1] = '{left_shoulder, right_shoulder} (y_displacement) {diff} - |spectral| rel_pwr_2_to_4'
$data[24,2] = 0.913
$data[24,3] = 0.01
$data[25,0] = 127
$data[25,1] = '{middle_shoulder, middle_hip} (x_displacement) {diff} - |spectral| rel_pwr_2_to_4'
$data[25,2] = 0.913
$data[25,3] = 0.058
$data[26,0] = 132
$data[26,1] = '{nose, middle_shoulder, left_shoulder} (angle) - |spectral| half_pwr_freq'
$data[26,2] = 0.957
$data[26,3] = 0.004
$data[27,0] = 134
$data[27,1] = '{nose, middle_shoulder, left_shoulder} (angle) - |spectral| rel_pwr_2_to_4'
$data[27,2] = 1
$data[27,3] = 0
$data[28,0] = 147
$data[28,1] = '{wrist, shoulder} (distance) - |spectral| rel_pwr_0.5_to_1_min'
$data[28,2] = 0.87
$data[28,3] = 1
$data[29,0] = 151
$data[29,1] = '{wrist, shoulder} (distance) - |spectral| rel_pwr_4_to_6_min'
$data[29,2] = 0.826
$data[29,3] = 0.143
$data[30,0] = 153
$data[30,1] = '{wrist, shoulder} (distance) - |spectral| rel_pwr_6_to_128_min'
$data[30,2] = 0.913
$data[30,3] = 0.431
$data[31,0] = 167
$data[31,1] = '{wrist, shoulder} (distance) {diff} - |temporal| mean_max'
$data[31,2] = 0.957
$data[31,3] = 0.008999999999999999

$ws.Range("A2:D33").Value = $data

# New rows (9-33) need the same column-A style (bold border, centered/top)
# that the existing A2:A8 cells already carry; copy it down.
$ws.Range("A8").Copy()
$ws.Range("A9:A33").PasteSpecial(-4122)  # xlPasteFormats
